# Fix for new ifc: add a new material row (Blocs béton manufacturés) to
# Sheet1, reusing the same "Nom_Matériau" list already used by the
# "Béton" row (B2), then leave the selection where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Blocs béton manufacturés"
$ws.Range("B6").Value = $ws.Range("B2").Value2

$ws.Range("N9").Select()
